# Rajout de message de felicitation en fin de partie
#
# The deck's closing slide has a subtitle placeholder that still reads the
# game's placeholder title ("Titre du jeu: a trouver"). Replace it with the
# end-of-game reveal text ("The hard disk"), bumping the point size to 36pt
# (PowerPoint reflows the placeholder with Shrink-text-on-overflow once the
# new text is in place).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder ("Sous-titre 2") robustly by name instead
# of assuming a fixed shape index.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Sous-titre 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$textRange = $shape.TextFrame.TextRange
$textRange.Text = "The hard "
$textRange.Font.Size = 36

# Append the second word as its own run (mirrors the source run split).
[void]$shape.TextFrame.TextRange.InsertAfter("disk")
